$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.424.10"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "'3.443.32"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'574.82"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "'144.32"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("D7").Value = "'3.443.83"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.480"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "'7.60"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("D12").Value = "'0.390"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "'4.027.24"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "'28.44"
$ws.Range("E14").Value = "  +4.35%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "'3.436.92"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "'61.510.61"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "'6.37"
$ws.Range("E19").Value = "  +4.14%  "
$ws.Range("D20").Value = "'14.41"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").Value = "'9.36"
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").Value = "'399.15"
$ws.Range("E22").Value = "  +3.80%  "
$ws.Range("D23").Value = "'0.564"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").Value = "'0.0000123"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").Value = "'3.583.19"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").Value = "'7.59"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "'1.49"
$ws.Range("E31").Value = "  -7.35%  "
$ws.Range("D32").Value = "'8.24"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'23.99"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'7.04"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "'3.469.89"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'5.16"
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("D39").Value = "'1.55"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "'166.69"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'0.0790"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "'27.70"
$ws.Range("E42").Value = "  +5.13%  "
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").Value = "'4.53"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.73"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'2.622.85"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.15"
$ws.Range("E48").Value = "  -4.63%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'6.95"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'23.16"
$ws.Range("E50").Value = "  -3.78%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'2.40"
$ws.Range("E51").Value = "  +2.60%  "
